$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nRows = 24

$block1 = New-Object "object[,]" $nRows,2
$block1[0,0] = 0.283596088079122
$block1[0,1] = 0.04750807013593317
$block1[1,0] = 0.2489195046054533
$block1[1,1] = 0.04283765755388913
$block1[2,0] = 0.2275817205730561
$block1[2,1] = 0.03995203726985608
$block1[3,0] = 0.2188753903838005
$block1[3,1] = 0.03877166957762768
$block1[4,0] = 0.21742906395761
$block1[4,1] = 0.03857540330402287
$block1[5,0] = 0.2274643477257143
$block1[5,1] = 0.03993613636667703
$block1[6,0] = 0.2716495524542495
$block1[6,1] = 0.04590147830609226
$block1[7,0] = 0.3579069966267809
$block1[7,1] = 0.05745463186691779
$block1[8,0] = 0.4210180854538521
$block1[8,1] = 0.06585211208303576
$block1[9,0] = 0.4496671187364711
$block1[9,1] = 0.06965222341472099
$block1[10,0] = 0.4605065122975702
$block1[10,1] = 0.0710883036311003
$block1[11,0] = 0.4581724825241906
$block1[11,1] = 0.07077914993180912
$block1[12,0] = 0.4505590738706076
$block1[12,1] = 0.06977042991931626
$block1[13,0] = 0.4458944006121328
$block1[13,1] = 0.06915217486022129
$block1[14,0] = 0.4191445285899817
$block1[14,1] = 0.06560335842343079
$block1[15,0] = 0.4027183707169399
$block1[15,1] = 0.06342111530904049
$block1[16,0] = 0.3932648230400844
$block1[16,1] = 0.0621640737773248
$block1[17,0] = 0.3900630632513753
$block1[17,1] = 0.06173814160308666
$block1[18,0] = 0.4044675535159001
$block1[18,1] = 0.06365361318242435
$block1[19,0] = 0.4527955762858937
$block1[19,1] = 0.07006679578270791
$block1[20,0] = 0.484325785537095
$block1[20,1] = 0.07424101479601575
$block1[21,0] = 0.4675027687487727
$block1[21,1] = 0.0720147498663124
$block1[22,0] = 0.4036767788542193
$block1[22,1] = 0.06354850849135119
$block1[23,0] = 0.334616343092847
$block1[23,1] = 0.05434494560125813
$ws.Range("B2:C25").Value = $block1

$block2 = New-Object "object[,]" $nRows,7
$block2[0,0] = 0.7626029043482419
$block2[0,1] = 2.100123961954324
$block2[0,2] = 0.3384098322087894
$block2[0,3] = 0.5193485153305204
$block2[0,4] = 0.3860694594142551
$block2[0,5] = 0.03217290880991897
$block2[0,6] = 0.2896375603817773
$block2[1,0] = 0.7331989198385997
$block2[1,1] = 2.070717013394471
$block2[1,2] = 0.3413539821906468
$block2[1,3] = 0.5238190168285541
$block2[1,4] = 0.3909906325850141
$block2[1,5] = 0.03248577531753405
$block2[1,6] = 0.2528093121956942
$block2[2,0] = 0.715530698613378
$block2[2,1] = 2.053897623029272
$block2[2,2] = 0.3434120902637616
$block2[2,3] = 0.5267800823144739
$block2[2,4] = 0.3942256303655753
$block2[2,5] = 0.03269701065143948
$block2[2,6] = 0.2301202447741559
$block2[3,0] = 0.7084277361718279
$block2[3,1] = 2.047354286775615
$block2[3,2] = 0.3443136614866802
$block2[3,3] = 0.5280411194799655
$block2[3,4] = 0.3955975457973313
$block2[3,5] = 0.03278789822358164
$block2[3,6] = 0.2208556722882093
$block2[4,0] = 0.7072541501935063
$block2[4,1] = 2.046286529208857
$block2[4,2] = 0.3444671616866444
$block2[4,3] = 0.5282537988025169
$block2[4,4] = 0.3958285894659639
$block2[4,5] = 0.03280328027920909
$block2[4,6] = 0.2193161921047562
$block2[5,0] = 0.7154345129945625
$block2[5,1] = 2.053808119623113
$block2[5,2] = 0.343423994703663
$block2[5,3] = 0.526796868892653
$block2[5,4] = 0.3942439153979986
$block2[5,5] = 0.03269821693293729
$block2[5,6] = 0.2299953739826606
$block2[6,0] = 0.7523842559941158
$block2[6,1] = 2.089727624898828
$block2[6,2] = 0.339372954736362
$block2[6,3] = 0.5208450947764049
$block2[6,4] = 0.3877219769793498
$block2[6,5] = 0.03227681153245143
$block2[6,6] = 0.2769553876908617
$block2[7,0] = 0.8279125732416759
$block2[7,1] = 2.169995416441651
$block2[7,2] = 0.3334196736700363
$block2[7,3] = 0.5108877501922251
$block2[7,4] = 0.3766266874066524
$block2[7,5] = 0.03160244103136556
$block2[7,6] = 0.3684152265126386
$block2[8,0] = 0.8852921147820041
$block2[8,1] = 2.234995054008607
$block2[8,2] = 0.3302653314925621
$block2[8,3] = 0.5046155834660198
$block2[8,4] = 0.3695092668254478
$block2[8,5] = 0.03119991076813378
$block2[8,6] = 0.4352034404555241
$block2[9,0] = 0.9118098613916885
$block2[9,1] = 2.265881789864608
$block2[9,2] = 0.329096433263814
$block2[9,3] = 0.5019885024373565
$block2[9,4] = 0.3664962660664806
$block2[9,5] = 0.03103702634121497
$block2[9,6] = 0.465493941352463
$block2[10,0] = 0.9219113833402304
$block2[10,1] = 2.27776778674442
$block2[10,2] = 0.3286921575933945
$block2[10,3] = 0.5010261974513739
$block2[10,4] = 0.3653876723995353
$block2[10,5] = 0.03097825997558346
$block2[10,6] = 0.4769504238270486
$block2[11,0] = 0.9197331777848916
$block2[11,1] = 2.275199475396676
$block2[11,2] = 0.3287775178014556
$block2[11,3] = 0.5012320011187725
$block2[11,4] = 0.3656249877421907
$block2[11,5] = 0.03099078665457
$block2[11,6] = 0.4744836930516954
$block2[12,0] = 0.9126397202140026
$block2[12,1] = 2.266855851041925
$block2[12,2] = 0.3290624038480559
$block2[12,3] = 0.5019086814500611
$block2[12,4] = 0.3664044125075421
$block2[12,5] = 0.03103213315541709
$block2[12,6] = 0.4664367549459598
$block2[13,0] = 0.9083025661763315
$block2[13,1] = 2.261769875240333
$block2[13,2] = 0.3292419036222469
$block2[13,3] = 0.5023274016939041
$block2[13,4] = 0.3668860488271655
$block2[13,5] = 0.0310578387914564
$block2[13,6] = 0.4615059463880016
$block2[14,0] = 0.8835674853171014
$block2[14,1] = 2.233003083765453
$block2[14,2] = 0.3303470839918461
$block2[14,3] = 0.5047918196129544
$block2[14,4] = 0.3697106999455428
$block2[14,5] = 0.03121096317666172
$block2[14,6] = 0.4332219803228554
$block2[15,0] = 0.8684997901390688
$block2[15,1] = 2.2156934086061
$block2[15,2] = 0.331093288178586
$block2[15,3] = 0.5063615792845866
$block2[15,4] = 0.3715011281185703
$block2[15,5] = 0.03131008502610477
$block2[15,6] = 0.4158466937739718
$block2[16,0] = 0.8598723550136782
$block2[16,1] = 2.205861385280116
$block2[16,2] = 0.3315475185617558
$block2[16,3] = 0.5072857511361164
$block2[16,4] = 0.372552086811023
$block2[16,5] = 0.03136900097036488
$block2[16,6] = 0.405844288745584
$block2[17,0] = 0.8569579676305352
$block2[17,1] = 2.202553725998754
$block2[17,2] = 0.3317056094451516
$block2[17,3] = 0.507602316256957
$block2[17,4] = 0.3729115548130935
$block2[17,5] = 0.03138927568543082
$block2[17,6] = 0.4024561897742558
$block2[18,0] = 0.8700997242738424
$block2[18,1] = 2.217523212651187
$block2[18,2] = 0.3310112618273067
$block2[18,3] = 0.5061922725576835
$block2[18,4] = 0.3713083445244063
$block2[18,5] = 0.03129933627199932
$block2[18,6] = 0.4176972176102822
$block2[19,0] = 0.9147216163427743
$block2[19,1] = 2.269301420429883
$block2[19,2] = 0.3289776839667411
$block2[19,3] = 0.5017090416855225
$block2[19,4] = 0.366174597838107
$block2[19,5] = 0.03101990954456291
$block2[19,6] = 0.4688007175938935
$block2[20,0] = 0.944233355759593
$block2[20,1] = 2.304248274976942
$block2[20,2] = 0.3278722679414656
$block2[20,3] = 0.4989685017983021
$block2[20,4] = 0.3630080628757959
$block2[20,5] = 0.03085427816333031
$block2[20,6] = 0.5021186357054432
$block2[21,0] = 0.9284504480661298
$block2[21,1] = 2.285495101341127
$block2[21,2] = 0.328441750774239
$block2[21,3] = 0.5004138417742681
$block2[21,4] = 0.3646808242162543
$block2[21,5] = 0.03094112230176904
$block2[21,6] = 0.484343885363046
$block2[22,0] = 0.86937628446411
$block2[22,1] = 2.216695585903409
$block2[22,2] = 0.3310482673657731
$block2[22,3] = 0.5062687485453594
$block2[22,4] = 0.3713954346529196
$block2[22,5] = 0.03130418977425187
$block2[22,6] = 0.4168606366044116
$block2[23,0] = 0.8071497108568053
$block2[23,1] = 2.147224960260019
$block2[23,2] = 0.3348164556981601
$block2[23,3] = 0.5133981241133938
$block2[23,4] = 0.3794467305635809
$block2[23,5] = 0.03176857582461778
$block2[23,6] = 0.343742709091913
$ws.Range("E2:K25").Value = $block2

$block3 = New-Object "object[,]" $nRows,2
$block3[0,0] = 1.134550080575771
$block3[0,1] = 1.650157203335866
$block3[1,0] = 1.142952505837698
$block3[1,1] = 1.665554202443033
$block3[2,0] = 1.148494373019808
$block3[2,1] = 1.675984465236382
$block3[3,0] = 1.15084907533182
$block3[3,1] = 1.680480278157845
$block3[4,0] = 1.151245893996965
$block3[4,1] = 1.681241621752861
$block3[5,0] = 1.148525739118654
$block3[5,1] = 1.676044103926856
$block3[6,0] = 1.137367859351528
$block3[6,1] = 1.655263344282858
$block3[7,0] = 1.118519703032618
$block3[7,1] = 1.622266614021243
$block3[8,0] = 1.106514568480527
$block3[8,1] = 1.602760892564561
$block3[9,0] = 1.101452001800176
$block3[9,1] = 1.594917992197693
$block3[10,0] = 1.099592175008581
$block3[10,1] = 1.592096432244716
$block3[11,0] = 1.099990176457574
$block3[11,1] = 1.592697503201677
$block3[12,0] = 1.101297845571814
$block3[12,1] = 1.594682885862014
$block3[13,0] = 1.102106285371555
$block3[13,1] = 1.595918319351355
$block3[14,0] = 1.106853434937548
$block3[14,1] = 1.603294196048054
$block3[15,0] = 1.109867708828602
$block3[15,1] = 1.608083106572252
$block3[16,0] = 1.111638964036594
$block3[16,1] = 1.610934529939115
$block3[17,0] = 1.112245127526649
$block3[17,1] = 1.611916620871895
$block3[18,0] = 1.109542951051452
$block3[18,1] = 1.607563281429194
$block3[19,0] = 1.100912198012658
$block3[19,1] = 1.594095702263346
$block3[20,0] = 1.095605213194773
$block3[20,1] = 1.586158809575096
$block3[21,0] = 1.098407137099173
$block3[21,1] = 1.590315665335638
$block3[22,0] = 1.10968965470645
$block3[22,1] = 1.607797988457236
$block3[23,0] = 1.12329462246003
$block3[23,1] = 1.630361947177391
$ws.Range("N2:O25").Value = $block3
